$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.774.25'
$ws.Range('E2').Value = '  -0.52%  '

$ws.Range('D3').Value = '3.841.42'
$ws.Range('E3').Value = '  +2.61%  '

$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').Value = '601.45'
$ws.Range('E5').Value = '  -0.04%  '

$ws.Range('D6').Value = '161.60'
$ws.Range('E6').Value = '  -3.08%  '

$ws.Range('D7').Value = '3.838.33'
$ws.Range('E7').Value = '  +2.60%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  -1.59%  '

$ws.Range('E10').Value = '  -1.15%  '

$ws.Range('D11').Value = '6.30'
$ws.Range('E11').Value = '  -1.82%  '

$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -0.33%  '

$ws.Range('D13').Value = '36.76'
$ws.Range('E13').Value = '  -3.13%  '

$ws.Range('E14').Value = '  -2.15%  '

$ws.Range('D15').Value = '4.490.29'
$ws.Range('E15').Value = '  +2.73%  '

$ws.Range('D16').Value = '3.844.25'
$ws.Range('E16').Value = '  +2.79%  '

$ws.Range('D17').Value = '68.893.17'
$ws.Range('E17').Value = '  -0.31%  '

$ws.Range('D18').Value = '7.52'
$ws.Range('E18').Value = '  +2.20%  '

$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.113'
$ws.Range('E19').Value = '  -0.19%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '11.38'
$ws.Range('E20').Value = '  +2.63%  '

$ws.Range('D21').Value = '17.12'
$ws.Range('E21').Value = '  -1.44%  '

$ws.Range('D22').Value = '483.99'
$ws.Range('E22').Value = '  -2.02%  '

$ws.Range('D23').Value = '0.718'
$ws.Range('E23').Value = '  -1.17%  '

$ws.Range('E24').Value = '  +4.26%  '

$ws.Range('D25').Value = '83.90'
$ws.Range('E25').Value = '  -1.09%  '

$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  -2.44%  '

$ws.Range('D27').Value = '12.07'
$ws.Range('E27').Value = '  -1.92%  '

$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Value = '0.999'
$ws.Range('E28').Value = '  -0.15%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '9.98'
$ws.Range('E29').Value = '  -1.26%  '

$ws.Range('E30').Value = '  -0.99%  '

$ws.Range('D31').Value = '7.91'

$ws.Range('D32').Value = '3.996.03'
$ws.Range('E32').Value = '  +2.78%  '

$ws.Range('E33').Value = '  -4.12%  '

$ws.Range('D34').Value = '32.09'
$ws.Range('E34').Value = '  +1.67%  '

$ws.Range('D35').Value = '3.791.38'
$ws.Range('E35').Value = '  +2.99%  '

$ws.Range('E36').Value = '  -1.52%  '

$ws.Range('E37').Value = '  +1.15%  '

$ws.Range('E38').Value = '  +3.30%  '

$ws.Range('D39').Value = '5.89'
$ws.Range('E39').Value = '  -1.08%  '

$ws.Range('E40').Value = '  +0.06%  '

$ws.Range('E41').Value = '  -1.80%  '

$ws.Range('D42').Value = '437.30'
$ws.Range('E42').Value = '  +1.57%  '

$ws.Range('D43').Value = '2.96'
$ws.Range('E43').Value = '  -1.06%  '

$ws.Range('D44').Value = '48.45'
$ws.Range('E44').Value = '  -0.73%  '

$ws.Range('E45').Value = '  -0.84%  '

$ws.Range('D47').Value = '8.37'
$ws.Range('E47').Value = '  -1.25%  '

$ws.Range('D48').Value = '26.42'
$ws.Range('E48').Value = '  +12.04%  '

$ws.Range('D49').Value = '143.04'
$ws.Range('E49').Value = '  +1.32%  '

$ws.Range('D50').Value = '2.826.47'
$ws.Range('E50').Value = '  +1.39%  '

$ws.Range('D51').Value = '0.0359'
$ws.Range('E51').Value = '  +1.99%  '
